# "Generate Report for Handoff"
# A new handoff run replaced the old GUID-named source file
# (c6159970-c0da-4760-9c8e-2a9162e7e16e) with a fresh one
# (eaa27c4c-c4d6-487b-bdcf-c2af93753bc9), refreshed the handoff xliff
# hashes/timestamps, and cleared out the stale handback (target) info
# since the new handoff hasn't been handed back yet.

$wb = $excel.ActiveWorkbook

$newGuid = "eaa27c4c-c4d6-487b-bdcf-c2af93753bc9"
$newHash = "4dd6dc82f6b594ff1d4b4256c62615e0076399f9"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value2 = "$newGuid.md"
$ws.Range("B2").Value2 = "e2e\$newGuid.md"
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}
$ws.Range("G2").Value2 = "2016-08-23 13:01:29"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$toRemove = @()
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$I$2') {
        $toRemove += $hl
    }
}
foreach ($hl in $toRemove) {
    $hl.Delete()
}

$ws.Range("A2").Value2 = "$newGuid.md"
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

$ws.Range("G2").Value2 = "$newGuid.$newHash.zh-cn.xlf"
$ws.Range("H2").Value2 = "2016-08-23 13:01:24"
$ws.Range("I2").Value2 = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value2 = ""
$ws.Range("K2").Value2 = "0001-01-01 00:00:00"

# Columns I/J shrink now that the target/handback-file text is gone
# (closest reproducible widths to the recorded 18.65 / 21.71 char units)
$ws.Columns.Item(9).ColumnWidth = 17.75
$ws.Columns.Item(10).ColumnWidth = 20.75

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$toRemove = @()
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$I$2') {
        $toRemove += $hl
    }
}
foreach ($hl in $toRemove) {
    $hl.Delete()
}

$ws.Range("A2").Value2 = "$newGuid.md"
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

$ws.Range("G2").Value2 = "$newGuid.$newHash.de-de.xlf"
$ws.Range("H2").Value2 = "2016-08-23 13:01:29"
$ws.Range("I2").Value2 = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value2 = ""
$ws.Range("K2").Value2 = "0001-01-01 00:00:00"

# Columns I/J shrink now that the target/handback-file text is gone
# (closest reproducible widths to the recorded 18.65 / 21.71 char units)
$ws.Columns.Item(9).ColumnWidth = 17.75
$ws.Columns.Item(10).ColumnWidth = 20.75
